$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBSSportsEmail")

# New column B: Language selector text (written in this order so the
# shared-string table lines up with the recorded workbook state)
$ws.Range("B1").Value = "Language"
$ws.Range("B3").Value = "Português (Brasil)"
$ws.Range("B4").Value = "Français"
$ws.Range("B2").Value = "Español"
$ws.Range("B5").Value = "日本"

# New column C: Footer text
$ws.Range("C2").Value = "©2021 ViacomCBS - Todos los derechos reservados"
$ws.Range("C4").Value = "©2021 ViacomCBS - Tous droits réservés"
$ws.Range("C1").Value = "Footer"
$ws.Range("C3").Value = "©2021 ViacomCBS - Todos os direitos reservados"
$ws.Range("C5").Value = "©2021 ViacomCBS - All rights reserved"

# Column widths for the new columns (values chosen so the engine's
# pixel-quantized ColumnWidth lands on the closest achievable width to the
# recorded 36.6640625 / 39.1640625 character widths)
$ws.Columns.Item(2).ColumnWidth = 35.83
$ws.Columns.Item(3).ColumnWidth = 38.33

# Update selection to match the recorded state after editing
$ws.Range("C11").Select()
